$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while preserving the cell's
# original (default/General) number format and avoiding Excel's
# automatic conversion of numeric-looking strings (e.g. "5.420")
# into actual numbers (which would drop trailing zeros, etc.).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '27.184.61'
$ws.Range("E2").Value = '  +0.95%  '

Set-TextValue $ws.Range("D3") '1.857.94'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("E4").Value = '  +1.31%  '

Set-TextValue $ws.Range("D5") '1.018'
$ws.Range("E5").Value = '  +1.13%  '

Set-TextValue $ws.Range("D6") '311.37'
$ws.Range("E6").Value = '  +0.45%  '

Set-TextValue $ws.Range("D7") '0.4797'
$ws.Range("E7").Value = '  +1.90%  '

Set-TextValue $ws.Range("D8") '0.3714'

Set-TextValue $ws.Range("D9") '0.07298'
$ws.Range("E9").Value = '  +1.99%  '

Set-TextValue $ws.Range("D10") '0.9348'
$ws.Range("E10").Value = '  +0.85%  '

Set-TextValue $ws.Range("D11") '20.09'
$ws.Range("E11").Value = '  +2.52%  '

Set-TextValue $ws.Range("D12") '0.07869'
$ws.Range("E12").Value = '  +2.15%  '

Set-TextValue $ws.Range("D13") '1.881.68'
$ws.Range("E13").Value = '  +4.33%  '

Set-TextValue $ws.Range("D14") '5.420'
$ws.Range("E14").Value = '  +2.59%  '

$ws.Range("E15").Value = '  +1.97%  '

Set-TextValue $ws.Range("D16") '89.96'
$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("E17").Value = '  +1.05%  '

Set-TextValue $ws.Range("D18") '0.000008727'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("E19").Value = '  +1.11%  '

Set-TextValue $ws.Range("D20") '27.226.78'
$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("E21").Value = '  +1.71%  '

$ws.Range("E22").Value = '  +1.34%  '

Set-TextValue $ws.Range("D23") '10.67'
$ws.Range("E23").Value = '  +0.42%  '

Set-TextValue $ws.Range("D24") '1.955'
$ws.Range("E24").Value = '  +1.08%  '

Set-TextValue $ws.Range("D25") '153.73'
$ws.Range("E25").Value = '  +1.09%  '

Set-TextValue $ws.Range("D26") '18.47'
$ws.Range("E26").Value = '  +1.20%  '

Set-TextValue $ws.Range("D27") '1.994'
$ws.Range("E27").Value = '  -1.17%  '

Set-TextValue $ws.Range("D28") '115.53'
$ws.Range("E28").Value = '  +0.96%  '

Set-TextValue $ws.Range("D29") '4.935'
$ws.Range("E29").Value = '  +1.08%  '

Set-TextValue $ws.Range("D30") '0.08890'
$ws.Range("E30").Value = '  +0.35%  '

Set-TextValue $ws.Range("D31") '3.315'
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("E32").Value = '  +0.19%  '

Set-TextValue $ws.Range("D33") '4.581'
$ws.Range("E33").Value = '  +2.26%  '

Set-TextValue $ws.Range("D34") '0.7367'
$ws.Range("E34").Value = '  -1.65%  '

Set-TextValue $ws.Range("D35") '2.682'
$ws.Range("E35").Value = '  -3.76%  '

$ws.Range("E36").Value = '  +3.15%  '

Set-TextValue $ws.Range("D37") '0.02014'
$ws.Range("E37").Value = '  +3.68%  '

Set-TextValue $ws.Range("D38") '2.991'
$ws.Range("E38").Value = '  +1.04%  '

Set-TextValue $ws.Range("D39") '0.05246'
$ws.Range("E39").Value = '  +0.74%  '

Set-TextValue $ws.Range("D40") '0.5334'
$ws.Range("E40").Value = '  +2.05%  '

Set-TextValue $ws.Range("D41") '7.074'
$ws.Range("E41").Value = '  +1.35%  '

Set-TextValue $ws.Range("D42") '0.1524'
$ws.Range("E42").Value = '  +0.35%  '

Set-TextValue $ws.Range("D43") '8.335'
$ws.Range("E43").Value = '  +2.11%  '

Set-TextValue $ws.Range("D44") '10.65'
$ws.Range("E44").Value = '  +1.57%  '

$ws.Range("E45").Value = '  +1.24%  '

$ws.Range("E46").Value = '  +1.16%  '

Set-TextValue $ws.Range("D47") '102.43'
$ws.Range("E47").Value = '  +1.42%  '

Set-TextValue $ws.Range("D48") '1.630'
$ws.Range("E48").Value = '  +1.98%  '

Set-TextValue $ws.Range("D49") '66.36'
$ws.Range("E49").Value = '  +1.08%  '

Set-TextValue $ws.Range("D50") '0.06076'
$ws.Range("E50").Value = '  +0.69%  '

Set-TextValue $ws.Range("D51") '0.8973'
$ws.Range("E51").Value = '  +0.14%  '
